# Auto-committed on 2022/12/02 週五 17:19:16.40
#
# The underlying workbook "CdCashFlow.xlsx" documents a DB layout (DBD /
# DBS sheets). This edit:
#   1. Fills in the previously-blank "備註說明" (remark) column for three
#      of the DBD rows (SEQ 2, 3, 5 -> InterestIncome / PrincipalAmortizeAmt
#      / DuePaymentAmt) with the note "BS060維護".
#   2. Leaves the DBD sheet active with cell D15 selected (instead of the
#      DBS sheet being the active/selected tab with C5/B13 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("DBD")

# Row 10 (SEQ 2 - InterestIncome), Row 11 (SEQ 3 - PrincipalAmortizeAmt),
# Row 13 (SEQ 5 - DuePaymentAmt): set the remark/note column (G).
$ws.Range("G10").Value = "BS060維護"
$ws.Range("G11").Value = "BS060維護"
$ws.Range("G13").Value = "BS060維護"

# Make DBD the active sheet/tab and leave D15 as the selected cell.
$ws.Activate()
$ws.Range("D15").Select()
